# Commit: "Change schema to English"
# The Staff sheet's header row (row 1) is translated from Chinese labels
# to their English equivalents. Everything else (data rows, Partners
# sheet, hyperlinks, column widths, cell styles, etc.) is left untouched.
#
# Header mapping (old Chinese -> new English):
#   A1 序号                                               -> No
#   B1 姓名                                               -> Name
#   C1 职务类型                                           -> Title
#   D1 教育背景                                           -> Study
#   E1 是否是华人                                         -> Chinese
#   F1 所获荣誉或Title                                    -> Awards
#   G1 是否是某学术期刊的编辑或学术会议的组委会成员/专业委员会 -> Fellow
#   H1 研究方向                                           -> Research area
#   I1 有几个博士学生/博士后学生/硕士学生/RA等团队成员      -> Students
#   J1 近三年学术文章数量                                  -> Papers
#   K1 兼职/以往雇主                                      -> Past-employer
#   L1 E-mail                                             -> E-mail (unchanged)
#
# The header row was tall (58.5pt) to accommodate the wrapped Chinese
# text; with the shorter English labels it is resized to 31.5pt, and the
# previously selected cell (H20) moves to H10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Staff")

$ws.Range("A1").Value = "No"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Title"
$ws.Range("D1").Value = "Study"
$ws.Range("E1").Value = "Chinese"
$ws.Range("F1").Value = "Awards"
$ws.Range("G1").Value = "Fellow"
$ws.Range("H1").Value = "Research area"
$ws.Range("I1").Value = "Students"
$ws.Range("J1").Value = "Papers"
$ws.Range("K1").Value = "Past-employer"
$ws.Range("L1").Value = "E-mail"

$ws.Rows.Item(1).RowHeight = 31.5

$ws.Activate()
$ws.Range("H10").Select()
